{"js": "// Change \"with a total of 18 players\" -> \"with a total of 19 players\"\n// in the ABSTRACT paragraph that discusses the questionnaire results.\nconst body = context.document.body;\n\nconst results = body.search(\"with a total of 18 players\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nresults.load(\"items,text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Could not find \"with a total of 18 players\" in the document body.');\n}\n\nresults.items[0].insertText(\"with a total of 19 players\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Change \"with a total of 18 players\" -> \"with a total of 19 players\"\n# in the ABSTRACT paragraph that discusses the questionnaire results.\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"with a total of 18 players\"\n$find.Replacement.Text = \"with a total of 19 players\"\n\n$wdReplaceOne = 1\n$wdFindContinue = 1\n$found = $find.Execute(\n    $find.Text,             # FindText\n    $false,                 # MatchCase\n    $false,                 # MatchWholeWord\n    $false,                 # MatchWildcards\n    $false,                 # MatchSoundsLike\n    $false,                 # MatchAllWordForms\n    $true,                  # Forward\n    $wdFindContinue,        # Wrap\n    $false,                 # Format\n    $find.Replacement.Text, # ReplaceWith\n    $wdReplaceOne           # Replace\n)\n\nif (-not $found) {\n    throw 'Could not find \"with a total of 18 players\" in the document.'\n}\n"}
